# Updated cryptos list on Wed Sep 18 04:44:37 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.471.73"
$ws.Range("E2").Value = "  +3.80%  "
$ws.Range("D3").Value = "2.327.39"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.05"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.11"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "2.326.40"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.335"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.74"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "60.456.58"
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("D16").Value = "2.742.27"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "2.339.62"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.15"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.78"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.99"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.87"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("E28").Value = "  +4.17%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.66"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  +10.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.74"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.93"
$ws.Range("E33").Value = "  +2.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.38"
$ws.Range("E34").Value = "  +11.86%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.86"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.06"
$ws.Range("E39").Value = "  +3.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "327.39"
$ws.Range("E40").Value = "  +13.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.53"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.97"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.13"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0943"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.30"
$ws.Range("E46").Value = "  +6.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0497"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.561"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0213"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "0.0₆0216"
$ws.Range("E50").Value = "  +17.68%  "
$ws.Range("E51").Value = "  +0.78%  "
